$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '49.904.79'
$ws.Range("E2").Value = '  +3.92%  '
$ws.Range("D3").Value = '2.656.22'
$ws.Range("E3").Value = '  +6.47%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '113.39'
$ws.Range("E5").Value = '  +7.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '328.09'
$ws.Range("E6").Value = '  +3.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.530'
$ws.Range("E7").Value = '  +1.93%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.555'
$ws.Range("E9").Value = '  +3.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.26'
$ws.Range("E10").Value = '  +6.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.25'
$ws.Range("E11").Value = '  +2.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0825'
$ws.Range("E12").Value = '  +2.94%  '
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.35'
$ws.Range("E14").Value = '  +4.79%  '
$ws.Range("D15").Value = '3.006.42'
$ws.Range("E15").Value = '  +3.82%  '
$ws.Range("D16").Value = '2.648.08'
$ws.Range("E16").Value = '  +5.64%  '
$ws.Range("E17").Value = '  +5.76%  '
$ws.Range("D18").Value = '49.823.24'
$ws.Range("E18").Value = '  +3.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.15'
$ws.Range("E19").Value = '  +2.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.77'
$ws.Range("E20").Value = '  +2.42%  '
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("D22").Value = '0.0₃0960'
$ws.Range("E22").Value = '  +3.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.34'
$ws.Range("E23").Value = '  +2.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '276.99'
$ws.Range("E24").Value = '  +2.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.59'
$ws.Range("E25").Value = '  +3.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.89'
$ws.Range("E26").Value = '  +4.95%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").Value = '  +2.96%  '
$ws.Range("E29").Value = '  +0.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.23'
$ws.Range("E30").Value = '  +5.16%  '
$ws.Range("E31").Value = '  -0.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.28'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.56'
$ws.Range("E33").Value = '  +2.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.42'
$ws.Range("E34").Value = '  +3.27%  '
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0801'
$ws.Range("E36").Value = '  +3.53%  '
$ws.Range("E37").Value = '  +7.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.77'
$ws.Range("E38").Value = '  +4.48%  '
$ws.Range("E39").Value = '  +7.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.43'
$ws.Range("E40").Value = '  +4.84%  '
$ws.Range("E41").Value = '  +2.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.23'
$ws.Range("E42").Value = '  +0.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.32'
$ws.Range("E43").Value = '  +0.47%  '
$ws.Range("E44").Value = '  +4.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.35'
$ws.Range("E45").Value = '  +5.44%  '
$ws.Range("D46").Value = '2.071.82'
$ws.Range("E46").Value = '  +3.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.26'
$ws.Range("E47").Value = '  +12.56%  '
$ws.Range("E48").Value = '  +6.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.16'
$ws.Range("E49").Value = '  +3.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.41'
$ws.Range("E50").Value = '  +5.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '82.12'
$ws.Range("E51").Value = '  +4.79%  '
